# Implemented general instance-based VCP notebook
#
# - "Cluster-Based" sheet: drop the leftover side-table (L15:P15 header row +
#   K16:K19 data) that used to live next to the Manhattan-distance block.
# - "KNN(Instance-based)" sheet: populate it with the RMSE / microseconds
#   results for the plain instance-based (instance-to-instance) KNN approach,
#   mirroring the layout already used on the "Cluster-Based" sheet.
# - Selection / active-tab bookkeeping: "KNN(Instance-based)" becomes the
#   selected tab (it was "Cluster-Based" before).

$wb = $excel.ActiveWorkbook

$wsCluster = $wb.Worksheets.Item("Cluster-Based")
$wsKnn = $wb.Worksheets.Item("KNN(Instance-based)")

# ---------------------------------------------------------------------------
# 1. Cluster-Based: remove the stray K/L..P side table (rows 15-19) that is
#    not part of the published results anymore.
# ---------------------------------------------------------------------------
[void]$wsCluster.Range("L15:P15").ClearContents()
[void]$wsCluster.Range("K16:K19").ClearContents()

# ---------------------------------------------------------------------------
# 2. KNN(Instance-based): build the "RMSE" block (rows 1-9) and the
#    "MICROSEGUNDOS" block (rows 10-17), same shape as on Cluster-Based but
#    with its own column set (B..F, with a couple of extra formatted columns
#    reserved up to I) and values for the instance-to-instance KNN run.
# ---------------------------------------------------------------------------

# -- Title row (merged A1:I1, bold + centered) ------------------------------
$r1 = $wsKnn.Range("A1:I1")
$r1.Font.Bold = $true
$r1.HorizontalAlignment = -4108
$wsKnn.Range("A1").Value = "BASADO EN KNN TRADICIONAL (INSTANCIA A INSTANCIA)"

# -- RMSE section header (merged A2:I2, centered) ---------------------------
$r2 = $wsKnn.Range("A2:I2")
$r2.HorizontalAlignment = -4108
$wsKnn.Range("A2").Value = "RMSE"

# -- RMSE column headers -----------------------------------------------------
$wsKnn.Range("B3").Value = "k=1"
$wsKnn.Range("C3").Value = "k=3"
$wsKnn.Range("D3").Value = "k=5"
$wsKnn.Range("E3").Value = "k=10"
$wsKnn.Range("F3").Value = "k=20"

# -- RMSE data rows -----------------------------------------------------------
$wsKnn.Range("A4").Value = "INV3-EUCLID"
$wsKnn.Range("B4").Value = 5.26
$wsKnn.Range("C4").Value = 4.31
$wsKnn.Range("D4").Value = 4.14
$wsKnn.Range("E4").Value = 4.09
$wsKnn.Range("F4").Value = 4.15

$wsKnn.Range("A5").Value = "INV3-COS"
$wsKnn.Range("B5").Value = 5.58
$wsKnn.Range("C5").Value = 4.36
$wsKnn.Range("D5").Value = 4.26
$wsKnn.Range("E5").Value = 4.07
$wsKnn.Range("F5").Value = 4.12

$wsKnn.Range("A6").Value = "HIST-EUCLID"
$wsKnn.Range("A7").Value = "HIST-COS"
$wsKnn.Range("A8").Value = "INV1-EUCLID"
$wsKnn.Range("A9").Value = "INV3-COS"

# -- MICROSEGUNDOS section header (merged A10:I10, centered) -----------------
$r10 = $wsKnn.Range("A10:I10")
$r10.HorizontalAlignment = -4108
$wsKnn.Range("A10").Value = "MICROSEGUNDOS"

# -- MICROSEGUNDOS column headers --------------------------------------------
$wsKnn.Range("B11").Value = "k=1"
$wsKnn.Range("C11").Value = "k=3"
$wsKnn.Range("D11").Value = "k=5"
$wsKnn.Range("E11").Value = "k=10"
$wsKnn.Range("F11").Value = "k=20"

# -- MICROSEGUNDOS data rows ---------------------------------------------------
$wsKnn.Range("A12").Value = "INV3-EUCLID"
$wsKnn.Range("B12").Value = 8915
$wsKnn.Range("C12").Value = 8764
$wsKnn.Range("D12").Value = 8976
$wsKnn.Range("E12").Value = 9000
$wsKnn.Range("F12").Value = 9390

$wsKnn.Range("A13").Value = "INV3-COS"
$wsKnn.Range("B13").Value = 11178
$wsKnn.Range("C13").Value = 11529
$wsKnn.Range("D13").Value = 11677
$wsKnn.Range("E13").Value = 11669
$wsKnn.Range("F13").Value = 11557

$wsKnn.Range("A14").Value = "HIST-EUCLID"
$wsKnn.Range("A15").Value = "HIST-COS"
$wsKnn.Range("A16").Value = "INV1-EUCLID"
$wsKnn.Range("A17").Value = "INV3-COS"

# -- Merge the title/section header rows -------------------------------------
[void]$wsKnn.Range("A2:I2").Merge()
[void]$wsKnn.Range("A1:I1").Merge()
[void]$wsKnn.Range("A10:I10").Merge()

# ---------------------------------------------------------------------------
# 3. Selection / active tab: KNN(Instance-based) is now the selected sheet.
# ---------------------------------------------------------------------------
[void]$wsCluster.Range("A2:F9").Select()
[void]$wsKnn.Activate()
[void]$wsKnn.Range("B14").Select()
